$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Barley / Low)
$ws.Range("D2").Value = 138
$ws.Range("E2").Value = 131
$ws.Range("F2").Value = 4.9
$ws.Range("G2").Value = 4.88
$ws.Range("H2").Value = 1200.8
$ws.Range("I2").Value = 1199.42
$ws.Range("J2").Value = 1.38

# Row 3 (Barley / Standard)
$ws.Range("D3").Value = 131
$ws.Range("E3").Value = 131
$ws.Range("F3").Value = 4.88
$ws.Range("G3").Value = 4.88
$ws.Range("H3").Value = 1152.27
$ws.Range("I3").Value = 1152.27
$ws.Range("J3").Value = 0

# Row 4 (Barley / High)
$ws.Range("D4").Value = 123.9
$ws.Range("E4").Value = 131
$ws.Range("F4").Value = 4.85
$ws.Range("G4").Value = 4.88
$ws.Range("H4").Value = 1106.49
$ws.Range("I4").Value = 1105.13
$ws.Range("J4").Value = 1.36

# Row 5 (Wheat / Low)
$ws.Range("E5").Value = 88.5
$ws.Range("F5").Value = 3.99
$ws.Range("G5").Value = 3.98
$ws.Range("H5").Value = 1055.51
$ws.Range("I5").Value = 1054.62
$ws.Range("J5").Value = 0.9

# Row 6 (Wheat / Standard)
$ws.Range("D6").Value = 88.5
$ws.Range("E6").Value = 88.5
$ws.Range("F6").Value = 3.98
$ws.Range("G6").Value = 3.98
$ws.Range("H6").Value = 1022.76
$ws.Range("I6").Value = 1022.76
$ws.Range("J6").Value = 0

# Row 7 (Wheat / High)
$ws.Range("D7").Value = 81.40000000000001
$ws.Range("E7").Value = 88.5
$ws.Range("F7").Value = 3.95
$ws.Range("G7").Value = 3.98
$ws.Range("H7").Value = 992.4
$ws.Range("I7").Value = 990.91
$ws.Range("J7").Value = 1.49

# Row 8 (Canola / Low)
$ws.Range("D8").Value = 164.6
$ws.Range("E8").Value = 146.9
$ws.Range("F8").Value = 1.71
$ws.Range("G8").Value = 1.68
$ws.Range("H8").Value = 802.85
$ws.Range("I8").Value = 799.7
$ws.Range("J8").Value = 3.16

# Row 9 (Canola / Standard)
$ws.Range("D9").Value = 146.9
$ws.Range("E9").Value = 146.9
$ws.Range("F9").Value = 1.68
$ws.Range("G9").Value = 1.68
$ws.Range("H9").Value = 746.8200000000001
$ws.Range("I9").Value = 746.8200000000001
$ws.Range("J9").Value = 0

# Row 10 (Canola / High)
$ws.Range("D10").Value = 129.2
$ws.Range("E10").Value = 146.9
$ws.Range("F10").Value = 1.65
$ws.Range("G10").Value = 1.68
$ws.Range("H10").Value = 697.3200000000001
$ws.Range("I10").Value = 693.9400000000001
$ws.Range("J10").Value = 3.38
